# Applies the per-coin "Price" (column D) and "Volume(1h)" (column E) updates
# from the Fri Nov 17 21:18:49 UTC 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'36.441.09"
$ws.Range('E2').Value = '  +1.32%  '
$ws.Range('D3').Value = "'1.944.70"
$ws.Range('E3').Value = '  -0.61%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = "'243.19"
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('E6').Value = '  -1.23%  '
$ws.Range('D7').Value = "'58.31"
$ws.Range('E7').Value = '  -2.82%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -0.93%  '
$ws.Range('D10').Value = "'55.81"
$ws.Range('E10').Value = '  -0.71%  '
$ws.Range('D11').Value = "'0.0833"
$ws.Range('E11').Value = '  +4.51%  '
$ws.Range('E12').Value = '  +1.25%  '
$ws.Range('D13').Value = "'21.63"
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').Value = "'0.821"
$ws.Range('E14').Value = '  -2.92%  '
$ws.Range('D15').Value = "'2.230.96"
$ws.Range('E15').Value = '  -0.40%  '
$ws.Range('D16').Value = "'13.61"
$ws.Range('E16').Value = '  -1.93%  '
$ws.Range('E17').Value = '  -2.31%  '
$ws.Range('D18').Value = "'1.944.49"
$ws.Range('E18').Value = '  -0.60%  '
$ws.Range('D19').Value = "'36.400.33"
$ws.Range('E19').Value = '  +1.53%  '
$ws.Range('D20').Value = "'69.70"
$ws.Range('E20').Value = '  -1.39%  '
$ws.Range('D21').Value = "'0.0₃0861"
$ws.Range('E21').Value = '  +1.52%  '
$ws.Range('D22').Value = "'229.50"
$ws.Range('E22').Value = '  -2.27%  '
$ws.Range('D23').Value = "'5.05"
$ws.Range('E23').Value = '  -2.24%  '
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('E25').Value = '  -3.70%  '
$ws.Range('E26').Value = '  +0.77%  '
$ws.Range('D27').Value = "'9.21"
$ws.Range('E27').Value = '  -4.85%  '
$ws.Range('D28').Value = "'162.18"
$ws.Range('E28').Value = '  +2.00%  '
$ws.Range('D29').Value = "'0.132"
$ws.Range('E29').Value = '  +1.47%  '
$ws.Range('D30').Value = "'19.45"
$ws.Range('E30').Value = '  -1.35%  '
$ws.Range('E31').Value = '  -0.78%  '
$ws.Range('E32').Value = '  +2.15%  '
$ws.Range('E33').Value = '  -2.81%  '
$ws.Range('D34').Value = "'0.0628"
$ws.Range('E34').Value = '  +2.32%  '
$ws.Range('E35').Value = '  -2.00%  '
$ws.Range('D36').Value = "'6.25"
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('D38').Value = "'1.77"
$ws.Range('E39').Value = '  -5.44%  '
$ws.Range('D40').Value = "'3.01"
$ws.Range('E40').Value = '  -0.59%  '
$ws.Range('D41').Value = "'0.0978"
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').Value = "'2.86"
$ws.Range('E42').Value = '  +0.99%  '
$ws.Range('E43').Value = '  -3.25%  '
$ws.Range('E44').Value = '  -0.60%  '
$ws.Range('D45').Value = "'16.01"
$ws.Range('E45').Value = '  +0.60%  '
$ws.Range('D46').Value = "'1.350.25"
$ws.Range('E46').Value = '  +1.53%  '
$ws.Range('E47').Value = '  -4.61%  '
$ws.Range('D48').Value = "'87.74"
$ws.Range('E48').Value = '  -4.26%  '
$ws.Range('D49').Value = "'7.11"
$ws.Range('E49').Value = '  -4.48%  '
$ws.Range('E50').Value = '  +1.12%  '
$ws.Range('D51').Value = "'45.42"
$ws.Range('E51').Value = '  +4.40%  '
